$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(706077, 699304, 79256, 498890, 79354, 394282, 459971, 514349, 514348, 223018)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
